# Update "想去人数" (number of interested attendees) figures for the
# 南宁 convention entries, reflecting refreshed counts as of 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2-6, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 10623
$wsExpo.Range("F3").Value = 236
$wsExpo.Range("F4").Value = 61
$wsExpo.Range("F5").Value = 677
$wsExpo.Range("F6").Value = 491

# Sheet "全部类型" (All types) - rows 2-5 and 7, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10623
$wsAll.Range("F3").Value = 236
$wsAll.Range("F4").Value = 61
$wsAll.Range("F5").Value = 677
$wsAll.Range("F7").Value = 491
